$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in columns B-E are text in the source workbook (inline strings),
# including numeric-looking Price values in column D. Force text format so Excel
# does not auto-convert them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.915.44"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.915.43"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.64"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4914"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2984"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06774"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.904.35"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.08"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07305"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.157"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.16"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6759"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.872.65"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007953"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.44"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.160.10"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.200"
$ws.Range("E22").Value = "  +6.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "207.53"
$ws.Range("E23").Value = "  +7.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.257"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.678"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.28"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.984"
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.333"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09188"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05190"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7526"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.714"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.737"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9255"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.100"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4523"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.58"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.009"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1406"
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.755"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.11"
$ws.Range("E47").Value = "  +15.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.033"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05953"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.95"
$ws.Range("E51").Value = "  +3.71%  "
